# Fix Training Data Issue (#48)
# The "Date" column (BF) for every data row (rows 2-31) was recorded as
# "6-3-2013-14" (a mangled month-day + season string). The correct value
# is the actual game date in ISO form: 2014-06-03.
#
# The date-like text must be written as literal text (not auto-converted
# to an Excel date serial number), matching the source file's use of
# plain date strings, so we force text entry with a leading apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.Value = "'2014-06-03"
